$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: fill in examination counts for "Dentar intraoral" row
$ws.Range("D5").Value = 650
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = "0df"
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 9
$ws.Range("J5").Value = 13
$ws.Range("K5").Value = 119
$ws.Range("L5").Value = 219
$ws.Range("M5").Value = 114
$ws.Range("N5").Value = 175

# Row 8 (TOTAL row): replace the SUM formulas for D:N with literal values
# (mirrors the app writing totals directly instead of relying on formulas)
$ws.Range("D8").Value = 650
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = "0gg"
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 9
$ws.Range("J8").Value = "1gggg"
$ws.Range("K8").Value = 119
$ws.Range("L8").Value = 219
$ws.Range("M8").Value = 114
$ws.Range("N8").Value = 175
